$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1025:1026 - this shifts the existing rows
# 1025..1152 down to 1027..1154 (matches the new dimension A1:R1154).
$ws.Rows("1025:1026").Insert()

# Row 1025: new "Primera" quality record dated 2023-07-17
$ws.Range("A1025").Value = 8
$ws.Range("B1025").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1025").Value = "Coquimbo"
$ws.Range("D1025").Value = "2023-07-17"
$ws.Range("E1025").Value = 4
$ws.Range("F1025").Value = 100112008
$ws.Range("G1025").Value = "Coliflor"
$ws.Range("H1025").Value = "Sin especificar"
$ws.Range("I1025").Value = "Primera"
$ws.Range("J1025").Value = 2000
$ws.Range("K1025").Value = 700
$ws.Range("L1025").Value = 800
$ws.Range("M1025").Value = 750
$ws.Range("N1025").Value = "$/unidad"
$ws.Range("O1025").Value = "Provincia del Elquí"
$ws.Range("P1025").Value = 750
$ws.Range("Q1025").Value = 1
$ws.Range("R1025").Value = "Hortaliza"

# Row 1026: new "Segunda" quality record, same date
$ws.Range("A1026").Value = 8
$ws.Range("B1026").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1026").Value = "Coquimbo"
$ws.Range("D1026").Value = "2023-07-17"
$ws.Range("E1026").Value = 4
$ws.Range("F1026").Value = 100112008
$ws.Range("G1026").Value = "Coliflor"
$ws.Range("H1026").Value = "Sin especificar"
$ws.Range("I1026").Value = "Segunda"
$ws.Range("J1026").Value = 1300
$ws.Range("K1026").Value = 500
$ws.Range("L1026").Value = 600
$ws.Range("M1026").Value = 550
$ws.Range("N1026").Value = "$/unidad"
$ws.Range("O1026").Value = "Provincia del Elquí"
$ws.Range("P1026").Value = 550
$ws.Range("Q1026").Value = 1
$ws.Range("R1026").Value = "Hortaliza"
